$wb = $excel.ActiveWorkbook

# Before: [1]=2021-Q3  [2]=总计
$wsTotalOrig = $wb.Worksheets.Item(2)

# --- Duplicate the "总计" sheet and insert the copy right before it, so the
#     final tab order becomes: 2021-Q3, 2022-Q1, 总计. Duplicating (instead of
#     Worksheets.Add) means the new sheet inherits the exact same styles
#     already present in the workbook (no new cellXfs entries needed).
#     NOTE: worksheet variables here track by *position*, not identity, so
#     after the insert shifts "总计" from index 2 to index 3 we must re-fetch
#     both sheets by name rather than keep using $wsTotalOrig. ---
$wsTotalOrig.Copy($wsTotalOrig)

$wsNew = $wb.Worksheets.Item("总计 (2)")
$wsNew.Name = "2022-Q1"
$wsTotal = $wb.Worksheets.Item("总计")

# $wsNew currently looks exactly like the old "总计" sheet:
#   B1:D1 = styled header cells, A2 = styled index cell, B2:D2 = data.
# Extend the styled header/index formatting across to column H by copying
# within the same sheet (cross-sheet Copy(Destination) does not transfer
# content reliably in this host, so every copy below targets a range on the
# SAME worksheet it was copied from).
$wsNew.Range("D1").Copy($wsNew.Range("E1:H1"))
$wsNew.Range("A2").Copy($wsNew.Range("A3:A4"))

# --- Header row ---
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# --- Data rows ---
$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "'001481"
$wsNew.Range("C2").Value = "华宝油气(QDII)美元"
$wsNew.Range("D2").Value = "'39.80"
$wsNew.Range("E2").Value = "'94.60"
$wsNew.Range("F2").Value = "'2.23"
$wsNew.Range("G2").Value = "'0.8875"
$wsNew.Range("H2").Value = 9

$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "'162411"
$wsNew.Range("C3").Value = "华宝油气(QDII)人民币A"
$wsNew.Range("D3").Value = "'39.80"
$wsNew.Range("E3").Value = "'94.60"
$wsNew.Range("F3").Value = "'2.23"
$wsNew.Range("G3").Value = "'0.8875"
$wsNew.Range("H3").Value = 9

$wsNew.Range("A4").Value = 2
$wsNew.Range("B4").Value = "'007844"
$wsNew.Range("C4").Value = "华宝油气(QDII)人民币C"
$wsNew.Range("D4").Value = "'12.98"
$wsNew.Range("E4").Value = "'94.60"
$wsNew.Range("F4").Value = "'2.23"
$wsNew.Range("G4").Value = "'0.2895"
$wsNew.Range("H4").Value = 9

# The leading apostrophes above force Excel to keep the numeric-looking
# strings as literal text (otherwise "001481"/"39.80" etc. get silently
# reinterpreted as numbers and lose their leading zeros / formatting), but
# they also leave a "quote prefix" flag on the cell's style. Strip that back
# off by pasting-in just the (plain, unstyled) format from C2 - a cell that
# already holds ordinary text with the default style - onto every cell that
# got the apostrophe treatment. This is a formats-only paste so the text
# values themselves are left untouched.
$wsNew.Range("C2").Copy()
$wsNew.Range("B2:B4").PasteSpecial(-4122)
$wsNew.Range("C2").Copy()
$wsNew.Range("D2:G4").PasteSpecial(-4122)

# --- Update the "总计" sheet: add the 2022-Q1 summary as the new row 2 and
#     push the existing 2021-Q3 row down to row 3. Values are written
#     directly (no Rows.Insert, which would fabricate a brand-new style). ---
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("B3").Value = "2021-Q3"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 1.84

$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 2.06
